$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at position A, shifting all existing columns (A:V) to (B:W).
$ws.Columns.Item(1).Insert()

# Row 2 (visible header row): new column A header = "Match ID", bold, no border (new style).
$ws.Range("A2").Value = "Match ID"
$ws.Range("A2").Font.Bold = $true

# Row 3 (hidden separator row): blank cell in column A, same bold/no-border style.
$ws.Range("A3").Font.Bold = $true

# Rows 4-19 (data rows): Match ID = 14 for every player row, same bold/no-border style.
$ws.Range("A4:A19").Value = 14
$ws.Range("A4:A19").Font.Bold = $true

# Row 20 (hidden summary row): Match ID = 14, default (unstyled) cell.
$ws.Cells.Item(20, 1).Value = 14
$ws.Rows.Item(20).AutoFit()

# Restore the worksheet selection to match the new layout.
$ws.Range("A2:A19").Select()
